# 2016-07-26  Create Slide Proxy
#
# - Add the "Tiến Độ" (progress) percentage values for the Slide row (F5)
#   and Product row (F7), formatted as a percentage (0%).
# - Move the active selection to J14 (was J9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slide row: 80% progress
$ws.Range("F5").Value = 0.8
$ws.Range("F5").NumberFormat = "0%"

# Product row: 30% progress
$ws.Range("F7").Value = 0.3
$ws.Range("F7").NumberFormat = "0%"

# Update the sheet's current selection/active cell
[void]$ws.Range("J14").Select()
